# The authored change swaps the two theme parts in this deck:
#   ppt/theme/theme1.xml  ("Office Theme" / "Office" color scheme)
#   ppt/theme/theme2.xml  ("Integral" / "Red Violet" color scheme)
# so that, afterwards, theme1.xml carries the Integral/Red-Violet colors
# and theme2.xml carries the Office/Office-Theme colors (font scheme and
# format scheme are identical between the two parts already, so only the
# 12 theme colors actually move).
#
# theme2.xml is the presentation's live design - it's the theme bound to
# the slide master (and therefore to every slide) via
# ppt/slideMasters/_rels/slideMaster1.xml.rels and ppt/_rels/presentation.xml.rels.
# The PowerPoint object model exposes exactly that theme's 12-color
# scheme through Slide.ThemeColorScheme (ThemeColor.RGB is read/write),
# so we repaint each of the 12 slots with the Office-theme RGB values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office theme colors (RGB() packs as r + g*256 + b*65536)
$tcs.Colors(1).RGB  = 0         # dk1      000000
$tcs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388   # dk2      44546A
$tcs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407     # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308  # accent5  4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink 954F72
